# Creating DataProvider Class attached to dataProvideClass annotation in testNG
$wb = $excel.ActiveWorkbook

# --- Fix the shared email string on RegistrationData (G1) ---
$regSheet = $wb.Worksheets.Item("RegistrationData")
$regSheet.Range("G1").Value = "ahmed.medhat@testautomation.com"

# --- Add the new LoginData sheet after RegistrationData ---
$loginSheet = $wb.Worksheets.Add($null, $regSheet)
$loginSheet.Name = "LoginData"

$loginSheet.Range("A1").Formula = "=RegistrationData!G1"
$loginSheet.Range("B1").Formula = "=CONCATENATE(""'"",RegistrationData!H1)"

$loginSheet.Range("A1").ColumnWidth = 41.28515625
$loginSheet.Range("B1").ColumnWidth = 13.140625

$loginSheet.Range("B1").Select() | Out-Null

# --- Adjust sheet selection state on RegistrationData ---
$regSheet.Range("C38").Select() | Out-Null

$loginSheet.Activate() | Out-Null
